# Open the workbook and select the forecast sheet that needs new columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KWH-2014-Forecast")

# --- Header row -------------------------------------------------------
# Column C's header changes from "KWH" to "Forecast_KWH", and four new
# columns (Lo80, Hi80, Lo95, Hi95) are added for the forecast intervals.
$ws.Range("C1").Value = "Forecast_KWH"
$ws.Range("D1").Value = "Lo80"
$ws.Range("E1").Value = "Hi80"
$ws.Range("F1").Value = "Lo95"
$ws.Range("G1").Value = "Hi95"

# --- Data rows ----------------------------------------------------------
# Updated forecast values (column C) plus the four new interval columns
# (D:G) for each of the 12 monthly rows (rows 2-13).
$data = @(
    @(9878216.82779281, 7851527.24526238, 11904906.4103232, 6778662.58713754, 12977771.0684481),
    @(7226198.35794253, 5735949.25104287, 8716447.4648422,  4947059.02711498, 9505337.68877008),
    @(7138953.97698387, 5659158.897073,   8618749.05689475, 4875802.70074758, 9402105.25322017),
    @(7063420.08973674, 5591859.39939538, 8534980.7800781,  4812862.22559937, 9313977.95387412),
    @(7248142.56535106, 5730517.16880919, 8765767.96189294, 4927134.81254841, 9569150.31815371),
    @(8797783.97566991, 6946534.87439008, 10649033.0769497, 5966542.79702977, 11629025.1543101),
    @(7599171.48630383, 5992263.83487319, 9206079.13773448, 5141618.30782714, 10056724.6647805),
    @(7693797.79159447, 6058946.16913086, 9328649.41405807, 5193507.99737626, 10194087.5858127),
    @(6801231.32952394, 5349058.85373009, 8253403.80531779, 4580325.1807442,  9022137.47830368),
    @(5116287.26016041, 4018649.07689258, 6213925.44342824, 3437594.5215294,  6794979.99879142),
    @(6997840.54727104, 5489418.48843338, 8506262.6061087,  4690908.08486074, 9304773.00968134),
    @(9152049.68401973, 7170006.94078882, 11134092.4272506, 6120776.89173008, 12183322.4763094)
)

$row = 2
foreach ($values in $data) {
    $ws.Cells.Item($row, 3).Value = $values[0]  # C - Forecast_KWH
    $ws.Cells.Item($row, 4).Value = $values[1]  # D - Lo80
    $ws.Cells.Item($row, 5).Value = $values[2]  # E - Hi80
    $ws.Cells.Item($row, 6).Value = $values[3]  # F - Lo95
    $ws.Cells.Item($row, 7).Value = $values[4]  # G - Hi95
    $row++
}
